# Generate Report for Handback
# Adds a new handback row (file c1370492-5503-4ef2-b839-a3332eb9aca0) to the
# "Overview", "zh-cn" and "de-de" worksheets, mirroring the layout of the
# two rows that already exist (c772ab87... / 63ef407f...).

$wb = $excel.ActiveWorkbook

$guid      = "c1370492-5503-4ef2-b839-a3332eb9aca0"
$mdName    = "$guid.md"
$xlfHash   = "7c8a99531402f758f54ffd8560be45dcace9f42e"
$statusOk  = "Handed back: in sync with en-US"
$includeTx = "Include"

# ---------------------------------------------------------------------------
# Helper: paint a cell the same way the existing "HyperLink" named style
# looks (blue underline) so new linked cells visually match A2/A3 etc.
# ---------------------------------------------------------------------------
function Set-LinkLook($cell) {
    $cell.Font.Underline = 2       # xlUnderlineStyleSingle
    $cell.Font.Color = 15570276    # RGB(0x64,0x95,0xED) == FF6495ED
}

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4, 2).Value = $statusOk
$wsOverview.Cells.Item(4, 3).Value = $statusOk

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
)
Set-LinkLook $wsOverview.Cells.Item(4, 1)

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Correspond Handoff File |
#                Correspond Handoff Datetime | Target File |
#                Correspond Handback File | Correspond Handback DateTime |
#                Handoff Reason | Dependency From
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhXlf = "$guid.$xlfHash.zh-cn.xlf"

$wsZh.Cells.Item(4, 2).Value = $statusOk
$wsZh.Cells.Item(4, 4).Value = "2016-02-24 08:47:07"
$wsZh.Cells.Item(4, 7).Value = "2016-02-24 08:47:48"
$wsZh.Cells.Item(4, 8).Value = $includeTx

$wsZh.Hyperlinks.Add($wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "", "", $mdName)
Set-LinkLook $wsZh.Cells.Item(4, 1)

$wsZh.Hyperlinks.Add($wsZh.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    "", "", $zhXlf)
Set-LinkLook $wsZh.Cells.Item(4, 3)

$wsZh.Hyperlinks.Add($wsZh.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "", "", $mdName)
Set-LinkLook $wsZh.Cells.Item(4, 5)

$wsZh.Hyperlinks.Add($wsZh.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    "", "", $zhXlf)
Set-LinkLook $wsZh.Cells.Item(4, 6)

# ---------------------------------------------------------------------------
# Sheet "de-de": same shape as "zh-cn"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$deXlf = "$guid.$xlfHash.de-de.xlf"

$wsDe.Cells.Item(4, 2).Value = $statusOk
$wsDe.Cells.Item(4, 4).Value = "2016-02-24 08:47:18"
$wsDe.Cells.Item(4, 7).Value = "2016-02-24 08:48:07"
$wsDe.Cells.Item(4, 8).Value = $includeTx

$wsDe.Hyperlinks.Add($wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "", "", $mdName)
Set-LinkLook $wsDe.Cells.Item(4, 1)

$wsDe.Hyperlinks.Add($wsDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    "", "", $deXlf)
Set-LinkLook $wsDe.Cells.Item(4, 3)

$wsDe.Hyperlinks.Add($wsDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "", "", $mdName)
Set-LinkLook $wsDe.Cells.Item(4, 5)

$wsDe.Hyperlinks.Add($wsDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    "", "", $deXlf)
Set-LinkLook $wsDe.Cells.Item(4, 6)
